$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the D16 cell text: the exam-task announcement is replaced with
# an "Orakel" (office hours) announcement for the same date (21.04).
$ws.Range("D16").Value = "21.04: Orakel med kursansvarlig i Aud Max og på  [Zoom](https://nhh.zoom.us/j/66065667678?pwd=ME1LK294VUw4SEt3eHI2V1ZuZm5MZz09)."

# Reflect the final selected/active cell as D16 (matches saved sheet view state).
$ws.Range("D16").Select()
